$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns F:H (Iteration - 4, 5, 6) entirely
$ws.Range("F1:H4").Delete() | Out-Null

# Update Mean row (row 2)
$ws.Range("B2").Value = 0.3219628369866471
$ws.Range("C2").Value = 0.08865399731597659
$ws.Range("D2").Value = 0.1300418914669286
$ws.Range("E2").Value = 0.0761832673032932

# Update Standard Deviation row (row 3)
$ws.Range("B3").Value = 1.67523729445211
$ws.Range("C3").Value = 0.5623973134614036
$ws.Range("D3").Value = 0.3392036008128436
$ws.Range("E3").Value = 0.2144067895675702

# Update Outlier row (row 4)
$ws.Range("B4").Value = 386
$ws.Range("C4").Value = 127
$ws.Range("D4").Value = 403
$ws.Range("E4").Value = 18
